$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values remain plain text when assigned (multi-dot price
# strings, percentage strings with surrounding spaces, etc.)
$textChanges = @{
    "D2" = "59.241.72"
    "E2" = "  -0.19%  "
    "D3" = "2.574.71"
    "E3" = "  -0.24%  "
    "E4" = "  -0.11%  "
    "E5" = "  -2.24%  "
    "E6" = "  -1.90%  "
    "E7" = "  -0.05%  "
    "E8" = "  -0.62%  "
    "D9" = "2.583.77"
    "E9" = "  -0.16%  "
    "E10" = "  +1.34%  "
    "E11" = "  +1.07%  "
    "E12" = "  +5.23%  "
    "E13" = "  +4.11%  "
    "D14" = "3.031.37"
    "D15" = "59.263.32"
    "E15" = "  -0.31%  "
    "E16" = "  +5.75%  "
    "E17" = "  +1.03%  "
    "D18" = "2.591.30"
    "E18" = "  +0.01%  "
    "E19" = "  +0.55%  "
    "E20" = "  +0.93%  "
    "E22" = "  +4.55%  "
    "E23" = "  -0.12%  "
    "E24" = "  +7.13%  "
    "E25" = "  -3.27%  "
    "E26" = "  -0.16%  "
    "E27" = "  -1.44%  "
    "E28" = "  +3.28%  "
    "D29" = "0.0₃0766"
    "E29" = "  -1.26%  "
    "E30" = "  -0.02%  "
    "E31" = "  -0.32%  "
    "E32" = "  +1.21%  "
    "E33" = "  -1.61%  "
    "E34" = "  +0.99%  "
    "E35" = "  +2.38%  "
    "E36" = "  +2.82%  "
    "E37" = "  +1.59%  "
    "E38" = "  +1.55%  "
    "E39" = "  -0.66%  "
    "E40" = "  -4.21%  "
    "E41" = "  +1.06%  "
    "E42" = "  -1.72%  "
    "E43" = "  +9.14%  "
    "E44" = "  -0.55%  "
    "E45" = "  -0.39%  "
    "E46" = "  +0.65%  "
    "E47" = "  +0.16%  "
    "E48" = "  -1.11%  "
    "E49" = "  +1.48%  "
    "D50" = "1.966.84"
    "E50" = "  +2.21%  "
    "E51" = "  +1.34%  "
}

# Price cells whose new values look like plain numbers (single decimal
# point). The source data models these as text (e.g. "552.06"), so force
# the Text number format before assigning, otherwise Excel would convert
# them to numeric values.
$numericLookingChanges = @{
    "D5" = "552.06"
    "D6" = "139.72"
    "D8" = "0.591"
    "D10" = "6.72"
    "D11" = "0.103"
    "D12" = "0.159"
    "D13" = "0.353"
    "D16" = "23.01"
    "D19" = "4.53"
    "D20" = "338.61"
    "D21" = "10.31"
    "D23" = "0.998"
    "D24" = "0.476"
    "D25" = "62.92"
    "D26" = "0.999"
    "D28" = "7.45"
    "D32" = "6.09"
    "D33" = "157.10"
    "D34" = "19.09"
    "D35" = "4.10"
    "D37" = "0.894"
    "D38" = "37.48"
    "D40" = "0.836"
    "D41" = "3.65"
    "D42" = "286.79"
    "D43" = "135.60"
    "D44" = "0.997"
    "D45" = "0.0969"
    "D46" = "0.594"
    "D47" = "10.67"
    "D48" = "0.0530"
    "D51" = "18.48"
}

foreach ($cell in $textChanges.Keys) {
    $ws.Range($cell).Value = $textChanges[$cell]
}

foreach ($cell in $numericLookingChanges.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $numericLookingChanges[$cell]
}
